# Updates of local parameter files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal TEXT value into a cell even when the text looks
# like a boolean/number (e.g. "True"/"False"), so Excel stores it as a
# shared string (t="s") instead of auto-converting it to a typed value.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Header row
$ws.Range("A1").Value = "model_id"
$ws.Range("B1").Value = "model"
$ws.Range("C1").Value = "local"
$ws.Range("D1").Value = "temperature"
$ws.Range("E1").Value = "system"

# Data rows: model_id, model, local, temperature, system
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "gemma2"
Set-TextValue $ws.Range("C2") "True"
$ws.Range("D2").Value = 0.7
$ws.Range("E2").Value = "All"

$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "gemma2:27b"
Set-TextValue $ws.Range("C3") "True"
$ws.Range("D3").Value = 0.7
$ws.Range("E3").Value = "Linux"

$ws.Range("A4").Value = 20
$ws.Range("B4").Value = "gpt-4o-mini"
Set-TextValue $ws.Range("C4") "False"
$ws.Range("D4").Value = 0.7
$ws.Range("E4").Value = "All"

$ws.Range("A5").Value = 25
$ws.Range("B5").Value = "gpt-4o"
Set-TextValue $ws.Range("C5") "False"
$ws.Range("D5").Value = 0.7
$ws.Range("E5").Value = "All"

$ws.Range("A6").Value = 30
$ws.Range("B6").Value = "llama3"
Set-TextValue $ws.Range("C6") "True"
$ws.Range("D6").Value = 0.7
$ws.Range("E6").Value = "All"

$ws.Range("A7").Value = 35
$ws.Range("B7").Value = "llama3:70b"
Set-TextValue $ws.Range("C7") "True"
$ws.Range("D7").Value = 0.7
$ws.Range("E7").Value = "Linux"

$ws.Range("A8").Value = 40
$ws.Range("B8").Value = "phi3"
Set-TextValue $ws.Range("C8") "True"
$ws.Range("D8").Value = 0.7
$ws.Range("E8").Value = "All"

$ws.Range("A9").Value = 45
$ws.Range("B9").Value = "phi3:medium"
Set-TextValue $ws.Range("C9") "True"
$ws.Range("D9").Value = 0.7
$ws.Range("E9").Value = "All"

# Column widths (best-fit-like custom widths, as observed in the target file)
$ws.Columns.Item(1).ColumnWidth = 7.833333333333333
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 4.666666666666667
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667

# Selection moves to F10 (just below/right of the table), as in the target file
$ws.Range("F10").Select() | Out-Null
